# Add a new "2022-Q3" sheet right after "总计" (i.e. before the current
# first-quarter sheet, "2022-Q2") and populate it with the fund holdings
# for that quarter, then update the "总计" summary sheet with a new
# leading row for 2022-Q3 (existing rows stay, shifted down by the
# worksheet itself since we are inserting a whole row).

$wb = $excel.ActiveWorkbook

$firstDataSheet = $wb.Worksheets.Item(2)   # "2022-Q2" before this edit
$q3 = $wb.Worksheets.Add($firstDataSheet)  # inserted immediately before it
$q3.Name = "2022-Q3"

# ---- Header row ----
# Columns B-H hold text headers; pre-format as text so short numeric-looking
# codes elsewhere in the column aren't auto-coerced to numbers later.
$q3.Range("B1:H1").NumberFormat = "@"
$q3.Range("B2:G9").NumberFormat = "@"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)  # headers start in column B
    $cell.Value = $headers[$i]
}
$headerRange = $q3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---- Data rows ----
# row: rowIndex, code, name, size, stockPosition, positionPct, marketValue, positionRank
# code/name/size/stockPosition/positionPct/marketValue (cols B-G) are text
# in the source data; rowIndex (A) and positionRank (H) are numeric.
$rows = @(
    @(0, "519700", "交银主题优选混合A",       "18.97", "70.23", "2.50", "0.4742", 8),
    @(1, "013884", "交银主题优选混合C",       "5.57",  "70.23", "2.50", "0.1392", 8),
    @(2, "001628", "招商体育文化休闲股票A",   "2.23",  "92.42", "5.01", "0.1117", 6),
    @(3, "516620", "国泰中证影视主题ETF",     "0.94",  "99.07", "7.14", "0.0671", 4),
    @(4, "159855", "银华中证影视主题ETF",     "0.84",  "96.84", "6.95", "0.0584", 4),
    @(5, "006048", "长城中证500指数增强A",    "2.12",  "94.78", "1.95", "0.0413", 5),
    @(6, "007413", "长城中证500指数增强C",    "1.06",  "94.78", "1.95", "0.0207", 5),
    @(7, "015395", "招商体育文化休闲股票C",   "0.25",  "92.42", "5.01", "0.0125", 6)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}

$aCol = $q3.Range("A2:A9")
$aCol.Font.Bold = $true
$aCol.HorizontalAlignment = -4108
$aCol.VerticalAlignment = -4160
$aCol.Borders.LineStyle = 1

# ---- Update the "总计" (summary) sheet ----
$summary = $wb.Worksheets.Item(1)

# Shift the existing 5 data rows (rows 2-6) down to rows 3-7 (bottom-up so
# we never overwrite a row we still need to read).
for ($row = 6; $row -ge 2; $row--) {
    $summary.Cells.Item($row + 1, 1).Value = $summary.Cells.Item($row, 1).Value()
    $summary.Cells.Item($row + 1, 2).Value = $summary.Cells.Item($row, 2).Value()
    $summary.Cells.Item($row + 1, 3).Value = $summary.Cells.Item($row, 3).Value()
    $summary.Cells.Item($row + 1, 4).Value = $summary.Cells.Item($row, 4).Value()
}

# New top row 2 with the 2022-Q3 aggregate numbers.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 8
$summary.Cells.Item(2, 4).Value = 0.93

# Re-number column A (row index) 0..5 and restore the bold/border style on
# the newly-populated A7 cell to match the rest of the column.
for ($row = 2; $row -le 7; $row++) {
    $summary.Cells.Item($row, 1).Value = $row - 2
}
$summary.Range("A2").Copy() | Out-Null
$summary.Range("A7").PasteSpecial(-4122) | Out-Null
